$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '81.429.27'
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('D3').Value = '3.165.65'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '210.40'
$ws.Range('E5').Value = '  +2.29%  '
$ws.Range('D6').Value = '621.27'
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('D7').Value = '0.279'
$ws.Range('E7').Value = '  +19.56%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '3.158.50'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  +9.71%  '
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '5.30'
$ws.Range('E14').Value = '  -4.04%  '
$ws.Range('D15').Value = '3.738.16'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').Value = '31.41'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '81.248.31'
$ws.Range('E17').Value = '  +2.57%  '
$ws.Range('D18').Value = '3.154.41'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '3.14'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('D20').Value = '13.91'
$ws.Range('E20').Value = '  -4.52%  '
$ws.Range('D21').Value = '431.28'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '8.93'
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('D23').Value = '5.07'
$ws.Range('E23').Value = '  +1.53%  '
$ws.Range('D24').Value = '7.22'
$ws.Range('E24').Value = '  +5.22%  '
$ws.Range('D25').Value = '5.25'
$ws.Range('E25').Value = '  +9.65%  '
$ws.Range('D26').Value = '3.317.76'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').Value = '76.21'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').Value = '10.76'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  +2.41%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '579.41'
$ws.Range('E32').Value = '  +10.29%  '
$ws.Range('D33').Value = '8.90'
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +8.66%  '
$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '1.99'
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('B37').Value = 'Cronos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D37').Value = '0.138'
$ws.Range('E37').Value = '  +14.58%  '
$ws.Range('D38').Value = '22.68'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '6.02'
$ws.Range('E40').Value = '  +10.19%  '
$ws.Range('D41').Value = '0.407'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '2.05'
$ws.Range('E42').Value = '  +13.87%  '
$ws.Range('D43').Value = '20.76'
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').Value = '3.01'
$ws.Range('E44').Value = '  +18.72%  '
$ws.Range('D45').Value = '159.12'
$ws.Range('E45').Value = '  -3.28%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '186.68'
$ws.Range('E47').Value = '  -3.05%  '
$ws.Range('D48').Value = '44.91'
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('D49').Value = '1.33'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').Value = '0.772'
$ws.Range('E50').Value = '  -4.81%  '
$ws.Range('D51').Value = '25.89'
$ws.Range('E51').Value = '  +0.21%  '

# Restore default (no explicit) formatting, matching original unstyled text cells
$ws.Range("D2:D51").ClearFormats()
